# Apply the CxSystem2 "post_syn_compartments" table update.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 (InhibitorySomaTarget / BC): Postsynaptic Compartment gains "nearestDendrite"
$ws.Range("C4").Value = "soma, nearestDendrite"

# Row 5 (InhibitoryDistalDendriteTarget / MC): Distribution weights updated
$ws.Range("E5").Value = ".33, .33, .33, 1, 1"

# Update the active selection to reflect where the author left the cursor
$ws.Range("E5").Select()
